$wb = $excel.ActiveWorkbook
$wsP = $wb.Worksheets.Item("Parameters")

# Mint style 7 (bold + vertical center) using A1 only
$wsP.Range("A1").Font.Bold = $true
$wsP.Range("A1").HorizontalAlignment = 1
$wsP.Range("A1").VerticalAlignment = -4108

# Mint style 8 (vertical center only) using A2 only
$wsP.Range("A2").HorizontalAlignment = 1
$wsP.Range("A2").VerticalAlignment = -4108

# Now copy those styles onto the remaining header / data rows
foreach ($r in @(4,7,10,13,16,19)) {
    $wsP.Range("A$r").Style = $wsP.Range("A1").Style
}
foreach ($r in @(3,5,6,8,9,11,12,14,15,17,18,20)) {
    $wsP.Range("A$r").Style = $wsP.Range("A2").Style
}

Write-Output "done"
